# "Fruta / hortaliza, semanal"
#
# A new weekly price record is inserted at the top of the data table
# (row 8, right after the header rows) for
#   Terminal La Palmera de La Serena - Arándano (blue).
# All the existing data rows (old rows 8-45) shift down by one row
# (becoming rows 9-46); the table grows from A1:T45 to A1:T46.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above row 8, shifting row 8 (and everything
# below it) down by one row. Formatting of the row above is carried
# down/through automatically (e.g. the date-formatted style in column D).
$ws.Rows.Item(8).Insert([Microsoft.Office.Interop.Excel.XlInsertShiftDirection]::xlShiftDown)

# Populate the newly inserted row 8 with the new weekly record.
$ws.Cells.Item(8, 1).Value  = 8
$ws.Cells.Item(8, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(8, 3).Value  = "Coquimbo"
$ws.Cells.Item(8, 4).Value  = 45037
$ws.Cells.Item(8, 5).Value  = 4
$ws.Cells.Item(8, 6).Value  = "Fruta"
$ws.Cells.Item(8, 7).Value  = 100101
$ws.Cells.Item(8, 8).Value  = "Berries"
$ws.Cells.Item(8, 9).Value  = 100101001
$ws.Cells.Item(8, 10).Value = "Arándano (blue)"
$ws.Cells.Item(8, 11).Value = "Sin especificar"
$ws.Cells.Item(8, 12).Value = "Primera"
$ws.Cells.Item(8, 13).Value = 160
$ws.Cells.Item(8, 14).Value = 10000
$ws.Cells.Item(8, 15).Value = 11000
$ws.Cells.Item(8, 16).Value = 10500
$ws.Cells.Item(8, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(8, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(8, 19).Value = 5250
$ws.Cells.Item(8, 20).Value = 2
